# tambahan spec baru konfirmasi pembayaran
# Rename "Sheet2" to "Transaksi" and fill it in with the new
# "Konfirmasi Pembayaran" master-table spec row (mirrors the existing
# Master sheet layout: Nama Tabel / Kode Tabel header + one data row).

$wb = $excel.ActiveWorkbook

$wsMaster = $wb.Worksheets.Item("Master")
$wsTrans  = $wb.Worksheets.Item("Sheet2")

$wsTrans.Name = "Transaksi"

$wsTrans.Range("A1").Value = "Nama Tabel"
$wsTrans.Range("B1").Value = "Kode Tabel"
$wsTrans.Range("A2").Value = "TRX001"
$wsTrans.Range("B2").Value = "Konfirmasi Pembayaran"

$wsTrans.Columns("A:B").AutoFit() | Out-Null

# Move the old selection on Master off of its previous spot ...
$wsMaster.Range("A1:B1").Select() | Out-Null

# ... and make the newly-populated Transaksi sheet the active tab,
# with the selection left below the data like in the source edit.
$wsTrans.Activate() | Out-Null
$wsTrans.Range("C11").Select() | Out-Null
